$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A1").Value = "basketball under pants"
$ws.Range("A2").Value = "softball gear for girls"
$ws.Range("A3").Value = "running capri"
$ws.Range("A4").Value = "softball compression sleeve"
$ws.Range("A5").Value = "youth softball compression sleeve"
$ws.Range("A6").Value = "running tights mens"
$ws.Range("A7").Value = "spandex men"
$ws.Range("A8").Value = "hockey kneepads"
$ws.Range("A9").Value = "padded leg sleeve"
$ws.Range("A10").Value = "mens basketball gear"
$ws.Range("A11").Value = "snowboarding padded shorts"
$ws.Range("A12").Value = "padded shorts snowboarding"
$ws.Range("A13").Value = "knee sleeve wrestling"
$ws.Range("A14").Value = "sleeve knee pads"
$ws.Range("A15").Value = "womens compression leggings"
$ws.Range("A16").Value = "airsoft knee pads"
$ws.Range("A17").Value = "mens compression tights 3 4"
$ws.Range("A18").Value = "basketball clothes for men"
$ws.Range("A19").Value = "men running tights"
$ws.Range("A20").Value = "knee pads nike"
$ws.Range("A21").Value = "knee pads mizuno"
$ws.Range("A22").Value = "knee pads bike"
$ws.Range("A23").Value = "yoga capri pants"
$ws.Range("A24").Value = "knee pads mtb"
$ws.Range("A25").Value = "knee pads skating"
$ws.Range("A26").Value = "mens workout tights"
$ws.Range("A27").Value = "mens basketball pants"
$ws.Range("A28").Value = "asics knee pads"
$ws.Range("A29").Value = "mens workout tights pants"
$ws.Range("A30").Value = "downhill knee pads"
$ws.Range("A31").Value = "men gym pants"
$ws.Range("A32").Value = "athletic capris"
$ws.Range("A33").Value = "valken knee pads"
$ws.Range("A34").Value = "woodland knee pads"
$ws.Range("A35").Value = "training tights men"
$ws.Range("A36").Value = "short tights for men"
$ws.Range("A37").Value = "ua compression pants"
$ws.Range("A38").Value = "men workout tights"
$ws.Range("A39").Value = "knee pads for exercise"
$ws.Range("A40").Value = "mens leggins"
$ws.Range("A41").Value = "nike kneepads"
$ws.Range("A42").Value = "youth football girdle with knee pads"
$ws.Range("A43").Value = "compression tights with pads"
$ws.Range("A44").Value = "compression leggings with knee pads"
$ws.Range("A45").Value = "basketball padded compression pants"
$ws.Range("A46").Value = "basketball knee pad pants"
$ws.Range("A47").Value = "padded compression pants men basketball"
$ws.Range("A48").Value = "padded tights men basketball"
$ws.Range("A49").Value = "tights with pads basketball"
$ws.Range("A50").Value = "mens leggings with knee pads"
$ws.Range("A51").Value = "leggings with knee pads women"
$ws.Range("A52").Value = "mtb knee pads men"
$ws.Range("A53").Value = "basketball padded knee sleeve"
$ws.Range("A54").Value = "compression pants women"
$ws.Range("A55").Value = "compression knee sleeve men basketball"
$ws.Range("A56").Value = "basketball sweat pants for men"
$ws.Range("A57").Value = "knee sleeve for wrestling"
$ws.Range("A58").Value = "leg sleeves for basketball youth"
$ws.Range("A59").Value = "training pants men"
$ws.Range("A60").Value = "compression knee sleeve men basketball"
$ws.Range("A61").Value = "basketball sweat pants for men"
$ws.Range("A62").Value = "knee sleeve for wrestling"
$ws.Range("A63").Value = "leg sleeves for basketball youth"
$ws.Range("A64").Value = "goalkeeper knee pads"
$ws.Range("A65").Value = "basketball calf sleeve"
$ws.Range("A66").Value = "compression knee sleeves with pads"
$ws.Range("A67").Value = "compression sleeve knee pads"
$ws.Range("A68").Value = "youth knee sleeve"
$ws.Range("A69").Value = "knee pad for scooter"
$ws.Range("A70").Value = "knee basketball"
$ws.Range("A71").Value = "knee pads for basketball youth"
$ws.Range("A72").Value = "mens compression knee"
$ws.Range("A73").Value = "knee pad sleeve basketball"
$ws.Range("A74").Value = "mens 3/4 compression pants"
$ws.Range("A75").Value = "youth compression knee pad sleeve"
$ws.Range("A76").Value = "mens basketball knee sleeves"
$ws.Range("A77").Value = "knee sleeve wrestling youth"
$ws.Range("A78").Value = "knee sleeves basketball men"
$ws.Range("A79").Value = "soccer compression pants"
$ws.Range("A80").Value = "leggings tight"
$ws.Range("A81").Value = "basketball leg sleeve youth padded"
$ws.Range("A82").Value = "knee pad construction"
$ws.Range("A83").Value = "youth basketball knee sleeve"
$ws.Range("A84").Value = "working knee pads for men"
$ws.Range("A85").Value = "cycling pants for men"
$ws.Range("A86").Value = "boys youth leggings"
$ws.Range("A87").Value = "compression running capris"
$ws.Range("A88").Value = "knee sleeve baseball"
$ws.Range("A89").Value = "compression knee sleeves for basketball"
$ws.Range("A90").Value = "volleyball kneepads"
$ws.Range("A91").Value = "compression knee sleeve with pad"
$ws.Range("A92").Value = "men capri shorts"
$ws.Range("A93").Value = "running compression pants"
$ws.Range("A94").Value = "mens work pants knee pads"
$ws.Range("A95").Value = "best knee pads"
$ws.Range("A96").Value = "compression pants sleeves"
$ws.Range("A97").Value = "mens compression running tights"
$ws.Range("A98").Value = "knee pads working"
$ws.Range("A99").Value = "basketball aids"
$ws.Range("A100").Value = "baseball youth pants"
